$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix/reword the "1 Korinther 13,12" bible verse text in C5
$ws.Range("C5").Value = '"Denn wir sehen jetzt durch einen Spiegel ein undeutliches Bild, dann aber, sehen wir von Angesicht zu Angesicht. Jetzt erkenne ich stückweise, dann aber werde ich erkennen, wie auch ich erkannt worden bin."  (1. Korinther 13,12)'

# 2. Add the new "Bibelstelle" column (D) with header + values, mirroring the
#    bold header style already used by B4/C4
$ws.Range("D4").Value = "Bibelstelle"
$ws.Range("D4").Font.Bold = $true

$ws.Range("D5").Value = "1 Kor 13,12"
$ws.Range("D6").Value = "Joh 16,22"
$ws.Range("D7").Value = "2 Tim 1,10"
$ws.Range("D8").Value = "1 Kor 15,27"
$ws.Range("D9").Value = "Jer 31,3"
$ws.Range("D10").Value = "Joh 10,27f."
$ws.Range("D11").Value = "Ps 118,6"
$ws.Range("D12").Value = "Ps 145,14"
$ws.Range("D13").Value = "Jes 43,1b"
$ws.Range("D14").Value = "Jes 54,10"
$ws.Range("D15").Value = "2 Kor 5,1"
$ws.Range("D16").Value = "1 Mos 24,56"
$ws.Range("D17").Value = "Ps 37,5"
$ws.Range("D18").Value = "5 Mos 4,31"

# 3. Set the width of the new column D to match the target (stored width 12)
$ws.Columns.Item(4).ColumnWidth = 11.2

# 4. Update the active selection to reflect where the user ended up (E7)
$ws.Range("E7").Select()
